$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Text)
    $origStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '27.350.06'
$ws.Range('E2').Value = '  +0.84%  '
Set-TextValue $ws.Range('D3') '1.824.60'
$ws.Range('E3').Value = '  -0.13%  '
Set-TextValue $ws.Range('D4') '1.000'
$ws.Range('E4').Value = '  -0.03%  '
Set-TextValue $ws.Range('D5') '314.09'
$ws.Range('E5').Value = '  +0.47%  '
Set-TextValue $ws.Range('D6') '0.9999'
$ws.Range('E6').Value = '  -0.07%  '
Set-TextValue $ws.Range('D7') '0.4463'
$ws.Range('E7').Value = '  -1.90%  '
Set-TextValue $ws.Range('D8') '0.3752'
$ws.Range('E8').Value = '  +0.44%  '
Set-TextValue $ws.Range('D9') '0.07479'
Set-TextValue $ws.Range('D10') '0.8856'
$ws.Range('E10').Value = '  +2.83%  '
Set-TextValue $ws.Range('D11') '21.03'
$ws.Range('E11').Value = '  +0.11%  '
Set-TextValue $ws.Range('D12') '1.834.90'
$ws.Range('E12').Value = '  +0.40%  '
Set-TextValue $ws.Range('D13') '6.754'
$ws.Range('E13').Value = '  +0.86%  '
Set-TextValue $ws.Range('D14') '5.423'
$ws.Range('E14').Value = '  +1.48%  '
Set-TextValue $ws.Range('D15') '93.80'
$ws.Range('E15').Value = '  +0.89%  '
Set-TextValue $ws.Range('D16') '0.07115'
$ws.Range('E16').Value = '  +0.49%  '
Set-TextValue $ws.Range('D17') '1.001'
$ws.Range('E17').Value = '  -0.11%  '
Set-TextValue $ws.Range('D18') '0.000008774'
$ws.Range('E18').Value = '  -0.73%  '
$ws.Range('E19').Value = '  -0.03%  '
Set-TextValue $ws.Range('D20') '15.15'
$ws.Range('E20').Value = '  +0.85%  '
Set-TextValue $ws.Range('D21') '27.352.63'
$ws.Range('E21').Value = '  +0.62%  '
Set-TextValue $ws.Range('D22') '5.416'
$ws.Range('E22').Value = '  +4.29%  '
$ws.Range('E23').Value = '  -0.72%  '
Set-TextValue $ws.Range('D24') '2.059.07'
$ws.Range('E24').Value = '  +0.40%  '
Set-TextValue $ws.Range('D25') '1.964'
$ws.Range('E25').Value = '  -2.08%  '
Set-TextValue $ws.Range('D26') '151.11'
$ws.Range('E26').Value = '  -0.57%  '
Set-TextValue $ws.Range('D27') '2.302'
$ws.Range('E27').Value = '  +3.05%  '
Set-TextValue $ws.Range('D28') '18.66'
$ws.Range('E28').Value = '  -0.06%  '
Set-TextValue $ws.Range('D29') '5.377'
$ws.Range('E29').Value = '  +2.07%  '
Set-TextValue $ws.Range('D30') '117.90'
$ws.Range('E30').Value = '  +0.26%  '
Set-TextValue $ws.Range('D31') '0.08889'
$ws.Range('E31').Value = '  +0.27%  '
Set-TextValue $ws.Range('D32') '0.7860'
$ws.Range('E32').Value = '  +3.39%  '
Set-TextValue $ws.Range('D33') '1.205'
$ws.Range('E33').Value = '  +0.78%  '
Set-TextValue $ws.Range('D34') '4.613'
$ws.Range('E34').Value = '  +3.06%  '
Set-TextValue $ws.Range('D35') '2.911'
$ws.Range('E35').Value = '  -2.48%  '
Set-TextValue $ws.Range('D36') '0.9993'
$ws.Range('E36').Value = '  -0.11%  '
Set-TextValue $ws.Range('D37') '1.113'
$ws.Range('E37').Value = '  +0.84%  '
Set-TextValue $ws.Range('D38') '0.01991'
$ws.Range('E38').Value = '  +1.01%  '
Set-TextValue $ws.Range('D39') '0.05295'
$ws.Range('E39').Value = '  +0.12%  '
Set-TextValue $ws.Range('D40') '7.309'
$ws.Range('E40').Value = '  +1.49%  '
Set-TextValue $ws.Range('D41') '0.5347'
$ws.Range('E41').Value = '  -0.64%  '
Set-TextValue $ws.Range('D42') '2.863'
$ws.Range('E42').Value = '  -1.02%  '
$ws.Range('E43').Value = '  +0.39%  '
Set-TextValue $ws.Range('D44') '2.303'
$ws.Range('E44').Value = '  +16.92%  '
Set-TextValue $ws.Range('D45') '8.662'
$ws.Range('E45').Value = '  +0.31%  '
Set-TextValue $ws.Range('D46') '0.5115'
$ws.Range('E46').Value = '  -1.81%  '
Set-TextValue $ws.Range('D47') '10.61'
$ws.Range('E47').Value = '  -1.28%  '
Set-TextValue $ws.Range('D48') '1.696'
$ws.Range('E48').Value = '  +1.06%  '
Set-TextValue $ws.Range('D49') '105.23'
$ws.Range('E49').Value = '  -0.99%  '
Set-TextValue $ws.Range('D50') '0.9993'
$ws.Range('E50').Value = '  -0.07%  '
Set-TextValue $ws.Range('D51') '0.06402'
$ws.Range('E51').Value = '  +0.72%  '
